$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant columns for all data rows (A,B,C,E,F,G,H,O,R stay fixed across rows)
$constB = "Agrícola del Norte S.A. de Arica"
$constC = "Arica y Parinacota"
$constG = "Jengibre"
$constH = "Sin especificar"
$constO = "Perú"
$constR = "Hortaliza"

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = $constB
$ws.Cells.Item(2, 3).Value = $constC
$ws.Cells.Item(2, 4).Value = 44320
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = 100114007
$ws.Cells.Item(2, 7).Value = $constG
$ws.Cells.Item(2, 8).Value = $constH
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 160
$ws.Cells.Item(2, 11).Value = 19000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 19500
$ws.Cells.Item(2, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(2, 15).Value = $constO
$ws.Cells.Item(2, 16).Value = 1500
$ws.Cells.Item(2, 17).Value = 13
$ws.Cells.Item(2, 18).Value = $constR

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = $constB
$ws.Cells.Item(3, 3).Value = $constC
$ws.Cells.Item(3, 4).Value = 44389
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = $constG
$ws.Cells.Item(3, 8).Value = $constH
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 13000
$ws.Cells.Item(3, 13).Value = 12500
$ws.Cells.Item(3, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value = $constO
$ws.Cells.Item(3, 16).Value = 962
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = $constR

# Row 4
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = $constB
$ws.Cells.Item(4, 3).Value = $constC
$ws.Cells.Item(4, 4).Value = 44918
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100114007
$ws.Cells.Item(4, 7).Value = $constG
$ws.Cells.Item(4, 8).Value = $constH
$ws.Cells.Item(4, 9).Value = "Segunda"
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 12000
$ws.Cells.Item(4, 12).Value = 13000
$ws.Cells.Item(4, 13).Value = 12750
$ws.Cells.Item(4, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(4, 15).Value = $constO
$ws.Cells.Item(4, 16).Value = 981
$ws.Cells.Item(4, 17).Value = 13
$ws.Cells.Item(4, 18).Value = $constR

# Row 5
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = $constB
$ws.Cells.Item(5, 3).Value = $constC
$ws.Cells.Item(5, 4).Value = 44764
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100114007
$ws.Cells.Item(5, 7).Value = $constG
$ws.Cells.Item(5, 8).Value = $constH
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 12000
$ws.Cells.Item(5, 12).Value = 13000
$ws.Cells.Item(5, 13).Value = 12500
$ws.Cells.Item(5, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(5, 15).Value = $constO
$ws.Cells.Item(5, 16).Value = 962
$ws.Cells.Item(5, 17).Value = 13
$ws.Cells.Item(5, 18).Value = $constR

# Row 6
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = $constB
$ws.Cells.Item(6, 3).Value = $constC
$ws.Cells.Item(6, 4).Value = 44988
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100114007
$ws.Cells.Item(6, 7).Value = $constG
$ws.Cells.Item(6, 8).Value = $constH
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 750
$ws.Cells.Item(6, 11).Value = 17000
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 17400
$ws.Cells.Item(6, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(6, 15).Value = $constO
$ws.Cells.Item(6, 16).Value = 1338
$ws.Cells.Item(6, 17).Value = 13
$ws.Cells.Item(6, 18).Value = $constR

# Row 7
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = $constB
$ws.Cells.Item(7, 3).Value = $constC
$ws.Cells.Item(7, 4).Value = 44379
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100114007
$ws.Cells.Item(7, 7).Value = $constG
$ws.Cells.Item(7, 8).Value = $constH
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 120
$ws.Cells.Item(7, 11).Value = 12000
$ws.Cells.Item(7, 12).Value = 13000
$ws.Cells.Item(7, 13).Value = 12667
$ws.Cells.Item(7, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(7, 15).Value = $constO
$ws.Cells.Item(7, 16).Value = 974
$ws.Cells.Item(7, 17).Value = 13
$ws.Cells.Item(7, 18).Value = $constR

# Row 8
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = $constB
$ws.Cells.Item(8, 3).Value = $constC
$ws.Cells.Item(8, 4).Value = 45092
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 100114007
$ws.Cells.Item(8, 7).Value = $constG
$ws.Cells.Item(8, 8).Value = $constH
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 11).Value = 13000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 13500
$ws.Cells.Item(8, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(8, 15).Value = $constO
$ws.Cells.Item(8, 16).Value = 1038
$ws.Cells.Item(8, 17).Value = 13
$ws.Cells.Item(8, 18).Value = $constR

# Row 9
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = $constB
$ws.Cells.Item(9, 3).Value = $constC
$ws.Cells.Item(9, 4).Value = 45049
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100114007
$ws.Cells.Item(9, 7).Value = $constG
$ws.Cells.Item(9, 8).Value = $constH
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 13500
$ws.Cells.Item(9, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(9, 15).Value = $constO
$ws.Cells.Item(9, 16).Value = 1038
$ws.Cells.Item(9, 17).Value = 13
$ws.Cells.Item(9, 18).Value = $constR

# Row 10
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = $constB
$ws.Cells.Item(10, 3).Value = $constC
$ws.Cells.Item(10, 4).Value = 44592
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = $constG
$ws.Cells.Item(10, 8).Value = $constH
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 120
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 13000
$ws.Cells.Item(10, 13).Value = 12500
$ws.Cells.Item(10, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = $constO
$ws.Cells.Item(10, 16).Value = 962
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = $constR

# Row 11
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = $constB
$ws.Cells.Item(11, 3).Value = $constC
$ws.Cells.Item(11, 4).Value = 44893
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = $constG
$ws.Cells.Item(11, 8).Value = $constH
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 900
$ws.Cells.Item(11, 11).Value = 13000
$ws.Cells.Item(11, 12).Value = 14000
$ws.Cells.Item(11, 13).Value = 13444
$ws.Cells.Item(11, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = $constO
$ws.Cells.Item(11, 16).Value = 1034
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = $constR

# Row 12
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = $constB
$ws.Cells.Item(12, 3).Value = $constC
$ws.Cells.Item(12, 4).Value = 44914
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100114007
$ws.Cells.Item(12, 7).Value = $constG
$ws.Cells.Item(12, 8).Value = $constH
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14400
$ws.Cells.Item(12, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(12, 15).Value = $constO
$ws.Cells.Item(12, 16).Value = 1108
$ws.Cells.Item(12, 17).Value = 13
$ws.Cells.Item(12, 18).Value = $constR

# Row 13
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = $constB
$ws.Cells.Item(13, 3).Value = $constC
$ws.Cells.Item(13, 4).Value = 44943
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100114007
$ws.Cells.Item(13, 7).Value = $constG
$ws.Cells.Item(13, 8).Value = $constH
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 350
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14429
$ws.Cells.Item(13, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(13, 15).Value = $constO
$ws.Cells.Item(13, 16).Value = 1110
$ws.Cells.Item(13, 17).Value = 13
$ws.Cells.Item(13, 18).Value = $constR

# Row 14
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = $constB
$ws.Cells.Item(14, 3).Value = $constC
$ws.Cells.Item(14, 4).Value = 45096
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100114007
$ws.Cells.Item(14, 7).Value = $constG
$ws.Cells.Item(14, 8).Value = $constH
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 750
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = 14600
$ws.Cells.Item(14, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(14, 15).Value = $constO
$ws.Cells.Item(14, 16).Value = 1123
$ws.Cells.Item(14, 17).Value = 13
$ws.Cells.Item(14, 18).Value = $constR

# Row 15
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = $constB
$ws.Cells.Item(15, 3).Value = $constC
$ws.Cells.Item(15, 4).Value = 44855
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100114007
$ws.Cells.Item(15, 7).Value = $constG
$ws.Cells.Item(15, 8).Value = $constH
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 500
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(15, 15).Value = $constO
$ws.Cells.Item(15, 16).Value = 769
$ws.Cells.Item(15, 17).Value = 13
$ws.Cells.Item(15, 18).Value = $constR

# Row 16
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = $constB
$ws.Cells.Item(16, 3).Value = $constC
$ws.Cells.Item(16, 4).Value = 44832
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = $constG
$ws.Cells.Item(16, 8).Value = $constH
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 14000
$ws.Cells.Item(16, 13).Value = 13500
$ws.Cells.Item(16, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(16, 15).Value = $constO
$ws.Cells.Item(16, 16).Value = 1038
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = $constR

# Row 17
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = $constB
$ws.Cells.Item(17, 3).Value = $constC
$ws.Cells.Item(17, 4).Value = 44406
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = 100114007
$ws.Cells.Item(17, 7).Value = $constG
$ws.Cells.Item(17, 8).Value = $constH
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 160
$ws.Cells.Item(17, 11).Value = 17000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 17500
$ws.Cells.Item(17, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(17, 15).Value = $constO
$ws.Cells.Item(17, 16).Value = 1346
$ws.Cells.Item(17, 17).Value = 13
$ws.Cells.Item(17, 18).Value = $constR

# Row 18
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = $constB
$ws.Cells.Item(18, 3).Value = $constC
$ws.Cells.Item(18, 4).Value = 44984
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100114007
$ws.Cells.Item(18, 7).Value = $constG
$ws.Cells.Item(18, 8).Value = $constH
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 16000
$ws.Cells.Item(18, 12).Value = 17000
$ws.Cells.Item(18, 13).Value = 16500
$ws.Cells.Item(18, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(18, 15).Value = $constO
$ws.Cells.Item(18, 16).Value = 1269
$ws.Cells.Item(18, 17).Value = 13
$ws.Cells.Item(18, 18).Value = $constR

# Row 19
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = $constB
$ws.Cells.Item(19, 3).Value = $constC
$ws.Cells.Item(19, 4).Value = 44890
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100114007
$ws.Cells.Item(19, 7).Value = $constG
$ws.Cells.Item(19, 8).Value = $constH
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14500
$ws.Cells.Item(19, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(19, 15).Value = $constO
$ws.Cells.Item(19, 16).Value = 1115
$ws.Cells.Item(19, 17).Value = 13
$ws.Cells.Item(19, 18).Value = $constR

# Row 20
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = $constB
$ws.Cells.Item(20, 3).Value = $constC
$ws.Cells.Item(20, 4).Value = 45028
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100114007
$ws.Cells.Item(20, 7).Value = $constG
$ws.Cells.Item(20, 8).Value = $constH
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 300
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 14500
$ws.Cells.Item(20, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(20, 15).Value = $constO
$ws.Cells.Item(20, 16).Value = 1115
$ws.Cells.Item(20, 17).Value = 13
$ws.Cells.Item(20, 18).Value = $constR

# Row 21
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = $constB
$ws.Cells.Item(21, 3).Value = $constC
$ws.Cells.Item(21, 4).Value = 44580
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = 100114007
$ws.Cells.Item(21, 7).Value = $constG
$ws.Cells.Item(21, 8).Value = $constH
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 160
$ws.Cells.Item(21, 11).Value = 11000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 11500
$ws.Cells.Item(21, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(21, 15).Value = $constO
$ws.Cells.Item(21, 16).Value = 885
$ws.Cells.Item(21, 17).Value = 13
$ws.Cells.Item(21, 18).Value = $constR

# Row 22
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = $constB
$ws.Cells.Item(22, 3).Value = $constC
$ws.Cells.Item(22, 4).Value = 45100
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100114007
$ws.Cells.Item(22, 7).Value = $constG
$ws.Cells.Item(22, 8).Value = $constH
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 15500
$ws.Cells.Item(22, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(22, 15).Value = $constO
$ws.Cells.Item(22, 16).Value = 1192
$ws.Cells.Item(22, 17).Value = 13
$ws.Cells.Item(22, 18).Value = $constR

# Row 23
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = $constB
$ws.Cells.Item(23, 3).Value = $constC
$ws.Cells.Item(23, 4).Value = 44972
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = 100114007
$ws.Cells.Item(23, 7).Value = $constG
$ws.Cells.Item(23, 8).Value = $constH
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 350
$ws.Cells.Item(23, 11).Value = 17000
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 13).Value = 17429
$ws.Cells.Item(23, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(23, 15).Value = $constO
$ws.Cells.Item(23, 16).Value = 1162
$ws.Cells.Item(23, 17).Value = 15
$ws.Cells.Item(23, 18).Value = $constR

# Row 24
$ws.Cells.Item(24, 1).Value = 1
$ws.Cells.Item(24, 2).Value = $constB
$ws.Cells.Item(24, 3).Value = $constC
$ws.Cells.Item(24, 4).Value = 44469
$ws.Cells.Item(24, 5).Value = 15
$ws.Cells.Item(24, 6).Value = 100114007
$ws.Cells.Item(24, 7).Value = $constG
$ws.Cells.Item(24, 8).Value = $constH
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 140
$ws.Cells.Item(24, 11).Value = 13000
$ws.Cells.Item(24, 12).Value = 14000
$ws.Cells.Item(24, 13).Value = 13500
$ws.Cells.Item(24, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(24, 15).Value = $constO
$ws.Cells.Item(24, 16).Value = 1038
$ws.Cells.Item(24, 17).Value = 13
$ws.Cells.Item(24, 18).Value = $constR

# Row 25
$ws.Cells.Item(25, 1).Value = 1
$ws.Cells.Item(25, 2).Value = $constB
$ws.Cells.Item(25, 3).Value = $constC
$ws.Cells.Item(25, 4).Value = 44910
$ws.Cells.Item(25, 5).Value = 15
$ws.Cells.Item(25, 6).Value = 100114007
$ws.Cells.Item(25, 7).Value = $constG
$ws.Cells.Item(25, 8).Value = $constH
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 50
$ws.Cells.Item(25, 11).Value = 14000
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14500
$ws.Cells.Item(25, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(25, 15).Value = $constO
$ws.Cells.Item(25, 16).Value = 1115
$ws.Cells.Item(25, 17).Value = 13
$ws.Cells.Item(25, 18).Value = $constR

# Row 26
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = $constB
$ws.Cells.Item(26, 3).Value = $constC
$ws.Cells.Item(26, 4).Value = 44159
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100114007
$ws.Cells.Item(26, 7).Value = $constG
$ws.Cells.Item(26, 8).Value = $constH
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 23000
$ws.Cells.Item(26, 12).Value = 24000
$ws.Cells.Item(26, 13).Value = 23500
$ws.Cells.Item(26, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(26, 15).Value = $constO
$ws.Cells.Item(26, 16).Value = 1808
$ws.Cells.Item(26, 17).Value = 13
$ws.Cells.Item(26, 18).Value = $constR

# Row 27
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = $constB
$ws.Cells.Item(27, 3).Value = $constC
$ws.Cells.Item(27, 4).Value = 45140
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100114007
$ws.Cells.Item(27, 7).Value = $constG
$ws.Cells.Item(27, 8).Value = $constH
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 200
$ws.Cells.Item(27, 11).Value = 16000
$ws.Cells.Item(27, 12).Value = 17000
$ws.Cells.Item(27, 13).Value = 16500
$ws.Cells.Item(27, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(27, 15).Value = $constO
$ws.Cells.Item(27, 16).Value = 1269
$ws.Cells.Item(27, 17).Value = 13
$ws.Cells.Item(27, 18).Value = $constR

# Row 28
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = $constB
$ws.Cells.Item(28, 3).Value = $constC
$ws.Cells.Item(28, 4).Value = 44397
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100114007
$ws.Cells.Item(28, 7).Value = $constG
$ws.Cells.Item(28, 8).Value = $constH
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 140
$ws.Cells.Item(28, 11).Value = 12500
$ws.Cells.Item(28, 12).Value = 13000
$ws.Cells.Item(28, 13).Value = 12750
$ws.Cells.Item(28, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(28, 15).Value = $constO
$ws.Cells.Item(28, 16).Value = 981
$ws.Cells.Item(28, 17).Value = 13
$ws.Cells.Item(28, 18).Value = $constR

# Row 29
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = $constB
$ws.Cells.Item(29, 3).Value = $constC
$ws.Cells.Item(29, 4).Value = 44616
$ws.Cells.Item(29, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100114007
$ws.Cells.Item(29, 7).Value = $constG
$ws.Cells.Item(29, 8).Value = $constH
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 120
$ws.Cells.Item(29, 11).Value = 19000
$ws.Cells.Item(29, 12).Value = 20000
$ws.Cells.Item(29, 13).Value = 19500
$ws.Cells.Item(29, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(29, 15).Value = $constO
$ws.Cells.Item(29, 16).Value = 1500
$ws.Cells.Item(29, 17).Value = 13
$ws.Cells.Item(29, 18).Value = $constR
